$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new instructional text in column E (rows 19-40) ---
$ws.Cells.Item(19, 5).Value = 'Для этого:'
$ws.Cells.Item(20, 5).Value = '1. Рассчитать средний уровень ряда.'
$ws.Cells.Item(21, 5).Value = '2. Рассчитать абсолютный прирост, темп роста и темп прироста для всех уровней и'
$ws.Cells.Item(22, 5).Value = 'средние значения этих показателей.'
$ws.Cells.Item(23, 5).Value = '3. Провести аналитическое выравнивание динамического ряда методом'
$ws.Cells.Item(24, 5).Value = 'наименьших квадратов, т.е. функцию y = f(x), заданную таблично,'
$ws.Cells.Item(25, 5).Value = 'аппроксимировать многочленом первой степени y = P1'
$ws.Cells.Item(27, 5).Value = '(x) = a1 + a2x;'
$ws.Cells.Item(29, 5).Value = '4. Рассчитать точность полученной динамической модели (уравнения) ряда,'
$ws.Cells.Item(30, 5).Value = 'сделать прогноз на два года.'
$ws.Cells.Item(31, 5).Value = '5. Определить колеблемость показателя и его устойчивость.'
$ws.Cells.Item(32, 5).Value = '6. Показать на графике исходный и выровненный динамические ряды.'
$ws.Cells.Item(33, 5).Value = '7. Средствами MS Excel построить точечную диаграмму по исходным табличным'
$ws.Cells.Item(34, 5).Value = 'данным, соответствующим своему варианту. Добавить линии тренда: линейный,'
$ws.Cells.Item(35, 5).Value = 'полиномиальный, экспоненциальный. Установить флаги «Показывать уравнение'
$ws.Cells.Item(36, 5).Value = 'на диаграмме», «Поместить на диаграмму величину достоверности'
$ws.Cells.Item(37, 5).Value = 'аппроксимации (R^2)». Для линейной линии тренда сравнить ее показатели с'
$ws.Cells.Item(38, 5).Value = 'рассчитанными аналитическим путем коэффициентами. Сделать вывод о том,'
$ws.Cells.Item(39, 5).Value = 'какая из линий тренда лучше аппроксимирует исходные данные.'
$ws.Cells.Item(40, 5).Value = '8. Сделать выводы по работе.'

# --- Adjust column widths (B, F, G) to match the widened layout ---
$ws.Columns.Item(2).ColumnWidth = 29.666666666666668
$ws.Columns.Item(6).ColumnWidth = 43.833333333333336
$ws.Columns.Item(7).ColumnWidth = 43.166666666666664

# --- Reposition/resize the two charts to make room for the new text ---
$co1 = $ws.ChartObjects(1)
$co1.Left = 1041.512645792323
$co1.Top = 203.57133858267716
$co1.Width = 655.0567688853346
$co1.Height = 271.9285826771653

$co2 = $ws.ChartObjects(2)
$co2.Left = 1041.7268190206694
$co2.Top = 474.8570866141732
$co2.Width = 755.3249387918306
$co2.Height = 521.1429133858268

# --- Update the active selection to reflect where the user was working ---
$ws.Range("D14").Select()
